# Update the "last modified" date shown in the date placeholders
# (2019/10/27 -> 2019/10/29) on every slide, and fix the title on the
# Overview slide that still read "1000 Overview" instead of "1100 Overview".

$p = $ppt.ActivePresentation

$oldDate = "2019/10/27"
$newDate = "2019/10/29"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)

        if (-not $sh.HasTextFrame) {
            continue
        }

        $tr = $sh.TextFrame.TextRange

        # Date placeholders (ppPlaceholderDate = 16): refresh the stale date.
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }

        # Title placeholder still carrying the old "1000 Overview" heading.
        if ($tr.Text -eq "1000 Overview") {
            $tr.Text = "1100 Overview"
        }
    }
}
